$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KO")

# Step 1: Insert two new columns before column D (shifts D:K -> F:M)
$ws.Range("D:E").Insert()

# Step 2: Copy cell formatting from column F (the old column D, now shifted) into new D:E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: Populate values for the two new columns (D, E) row by row
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43371
$ws.Range("D8").Value = 7058000
$ws.Range("E8").Value = 8245000
$ws.Range("D9").Value = 2721000
$ws.Range("E9").Value = 3059000
$ws.Range("D10").Value = 4337000
$ws.Range("E10").Value = 5186000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 667000
$ws.Range("E14").Value = 597000
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 5920000
$ws.Range("E17").Value = 6179000
$ws.Range("D18").Value = 1138000
$ws.Range("E18").Value = 2066000
$ws.Range("D20").Value = -114000
$ws.Range("E20").Value = 990000
$ws.Range("D21").Value = 1303000
$ws.Range("E21").Value = 3310000
$ws.Range("D22").Value = 237000
$ws.Range("E22").Value = 209000
$ws.Range("D23").Value = 787000
$ws.Range("E23").Value = 2847000
$ws.Range("D24").Value = -4000
$ws.Range("E24").Value = 653000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 791000
$ws.Range("E26").Value = 2194000
$ws.Range("D27").Value = 734000
$ws.Range("E27").Value = 2256000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 136000
$ws.Range("E29").Value = -376000
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 114000
$ws.Range("E32").Value = -990000
$ws.Range("D33").Value = 870000
$ws.Range("E33").Value = 1880000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 870000
$ws.Range("E35").Value = 1880000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43371
$ws.Range("D41").Value = 8926000
$ws.Range("E41").Value = 9065000
$ws.Range("D42").Value = 7038000
$ws.Range("E42").Value = 9782000
$ws.Range("D43").Value = 3396000
$ws.Range("E43").Value = 3702000
$ws.Range("D44").Value = 2766000
$ws.Range("E44").Value = 2627000
$ws.Range("D45").Value = 8508000
$ws.Range("E45").Value = 8237000
$ws.Range("D46").Value = 30634000
$ws.Range("E46").Value = 33413000
$ws.Range("D47").Value = 20274000
$ws.Range("E47").Value = 21950000
$ws.Range("D48").Value = 8232000
$ws.Range("E48").Value = 7404000
$ws.Range("D49").Value = 17270000
$ws.Range("E49").Value = 16855000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 6806000
$ws.Range("E52").Value = 7255000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 83216000
$ws.Range("E54").Value = 86877000
$ws.Range("D57").Value = 8764000
$ws.Range("E57").Value = 10253000
$ws.Range("D58").Value = 18191000
$ws.Range("E58").Value = 19314000
$ws.Range("D59").Value = 2268000
$ws.Range("E59").Value = 1863000
$ws.Range("D60").Value = 29223000
$ws.Range("E60").Value = 31430000
$ws.Range("D61").Value = 25364000
$ws.Range("E61").Value = 25523000
$ws.Range("D62").Value = 9571000
$ws.Range("E62").Value = 9746000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 66235000
$ws.Range("E66").Value = 68613000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 63234000
$ws.Range("E72").Value = 64028000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 16981000
$ws.Range("E76").Value = 18264000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43371
$ws.Range("D81").Value = 870000
$ws.Range("E81").Value = 1880000
$ws.Range("D83").Value = 279000
$ws.Range("E83").Value = 254000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 1937000
$ws.Range("E89").Value = 3004000
$ws.Range("D91").Value = -430000
$ws.Range("E91").Value = -305000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 1690000
$ws.Range("E94").Value = 1841000
$ws.Range("D96").Value = -3323000
$ws.Range("E96").Value = -1659000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -3835000
$ws.Range("E100").Value = -3535000
$ws.Range("D101").Value = -13000
$ws.Range("E101").Value = -140000
$ws.Range("D102").Value = -221000
$ws.Range("E102").Value = 1170000

# Step 4: Apply corrections to shifted cells (F:M) where the restated value differs from a pure shift
$ws.Range("H9").Value = 2688000
$ws.Range("I9").Value = 3394000
$ws.Range("H10").Value = 4824000
$ws.Range("I10").Value = 5684000
$ws.Range("F14").Value = 368000
$ws.Range("G14").Value = 589000
$ws.Range("H14").Value = 1201000
$ws.Range("I14").Value = 838000
$ws.Range("J14").Value = 795000
$ws.Range("F17").Value = 6341000
$ws.Range("G17").Value = 5869000
$ws.Range("H17").Value = 6975000
$ws.Range("I17").Value = 7069000
$ws.Range("J17").Value = 7387000
$ws.Range("F18").Value = 2586000
$ws.Range("G18").Value = 1757000
$ws.Range("H18").Value = 537000
$ws.Range("I18").Value = 2009000
$ws.Range("J18").Value = 2315000
$ws.Range("F20").Value = 530000
$ws.Range("G20").Value = 304000
$ws.Range("H20").Value = 609000
$ws.Range("I20").Value = -125000
$ws.Range("J20").Value = 540000
$ws.Range("H21").Value = 1480000
$ws.Range("I21").Value = 2181000
$ws.Range("H22").Value = 209000
$ws.Range("I22").Value = 210000
$ws.Range("I29").Value = 0
$ws.Range("F32").Value = -530000
$ws.Range("G32").Value = -304000
$ws.Range("H32").Value = -609000
$ws.Range("I32").Value = 125000
$ws.Range("J32").Value = -540000
$ws.Range("I91").Value = -362000
$ws.Range("J91").Value = -390000
$ws.Range("H94").Value = -2583000
$ws.Range("I94").Value = 1472000
$ws.Range("H102").Value = -6450000
